$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "new_new_name"

$values = @(
    "GT119",
    "GT120",
    "GT121",
    "GT122",
    "GT123",
    "GT124",
    "GT125",
    "GT126",
    "GT127",
    "GT128",
    "GT129",
    "GT130",
    "GT131",
    "GT132",
    "GT133",
    "GT134",
    "GT135"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

$ws.Range("D9").Select()
